$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text edit: "Ready for handoff" -> "In Translation" -------------------
# The same shared string is referenced from the Overview sheet (columns E
# and F, row 2) as well as the Status column (C) on both the zh-cn and
# de-de sheets, so every referencing cell needs to be updated.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width edit -----------------------------------------------------
# The "Status" columns shrink from ~17.22 characters wide to ~13.41 once the
# shorter "In Translation" text replaces "Ready for handoff". Reproduce this
# via the genuine ColumnWidth (character units) property.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
